$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: "ReciepeID" header + per-row numbering
$ws.Range("L1").Value = "ReciepeID"
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 2
$ws.Range("L4").Value = 3

# Row 4 (recipe 3) gets a name in column B
$ws.Range("B4").Value = "Desert"

# Update the view: select C4 (clears the stale topLeftCell/selection from T3)
$ws.Range("C4").Select()
